$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "51.801.45"
$ws.Range("E2").Value2 = "  -0.16%  "
$ws.Range("D3").Value2 = "2.959.79"
$ws.Range("E3").Value2 = "  +0.98%  "
$ws.Range("E4").Value2 = "  +0.05%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value2 = "351.99"
$c.Style = "Normal"
$ws.Range("E5").Value2 = "  -0.12%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value2 = "107.58"
$c.Style = "Normal"
$ws.Range("E6").Value2 = "  -4.17%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value2 = "0.556"
$c.Style = "Normal"
$ws.Range("E7").Value2 = "  -0.75%  "
$ws.Range("E8").Value2 = "  -0.01%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value2 = "0.608"
$c.Style = "Normal"
$ws.Range("E9").Value2 = "  -2.11%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value2 = "38.02"
$c.Style = "Normal"
$ws.Range("E10").Value2 = "  -3.59%  "
$ws.Range("E11").Value2 = "  +1.21%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value2 = "0.0850"
$c.Style = "Normal"
$ws.Range("E12").Value2 = "  -4.05%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value2 = "19.02"
$c.Style = "Normal"
$ws.Range("E13").Value2 = "  -5.54%  "
$ws.Range("D14").Value2 = "3.415.62"
$ws.Range("E14").Value2 = "  +0.82%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value2 = "7.62"
$c.Style = "Normal"
$ws.Range("E15").Value2 = "  -2.05%  "
$ws.Range("D16").Value2 = "2.937.81"
$ws.Range("E16").Value2 = "  +0.45%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value2 = "0.983"
$c.Style = "Normal"
$ws.Range("E17").Value2 = "  -0.20%  "
$ws.Range("D18").Value2 = "51.788.00"
$ws.Range("E18").Value2 = "  -0.29%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value2 = "3.39"
$c.Style = "Normal"
$ws.Range("E19").Value2 = "  +2.45%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value2 = "7.42"
$c.Style = "Normal"
$ws.Range("E20").Value2 = "  -2.48%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value2 = "13.46"
$c.Style = "Normal"
$ws.Range("E21").Value2 = "  -5.49%  "
$ws.Range("D22").Value2 = "0.0₃0968"
$ws.Range("E22").Value2 = "  -1.61%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value2 = "69.16"
$c.Style = "Normal"
$ws.Range("E23").Value2 = "  -2.91%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value2 = "261.95"
$c.Style = "Normal"
$ws.Range("E24").Value2 = "  -2.44%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value2 = "2.71"
$c.Style = "Normal"
$ws.Range("E25").Value2 = "  -2.35%  "
$ws.Range("E26").Value2 = "  -3.58%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value2 = "26.61"
$c.Style = "Normal"
$ws.Range("E27").Value2 = "  -1.35%  "
$ws.Range("E28").Value2 = "  +0.10%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value2 = "7.33"
$c.Style = "Normal"
$ws.Range("E29").Value2 = "  +1.45%  "
$ws.Range("E30").Value2 = "  +1.14%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value2 = "10.26"
$c.Style = "Normal"
$ws.Range("E31").Value2 = "  -3.02%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value2 = "6.07"
$c.Style = "Normal"
$ws.Range("E32").Value2 = "  -2.47%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value2 = "35.85"
$c.Style = "Normal"
$ws.Range("E33").Value2 = "  -3.66%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value2 = "2.16"
$c.Style = "Normal"
$ws.Range("E34").Value2 = "  -4.52%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value2 = "50.56"
$c.Style = "Normal"
$ws.Range("E35").Value2 = "  -4.42%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value2 = "0.0432"
$c.Style = "Normal"
$ws.Range("E36").Value2 = "  -4.67%  "
$ws.Range("E37").Value2 = "  -0.04%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value2 = "3.18"
$c.Style = "Normal"
$ws.Range("E38").Value2 = "  -4.99%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value2 = "17.75"
$c.Style = "Normal"
$ws.Range("E39").Value2 = "  -5.27%  "
$ws.Range("E40").Value2 = "  -4.80%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value2 = "2.67"
$c.Style = "Normal"
$ws.Range("E41").Value2 = "  -0.72%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value2 = "0.116"
$c.Style = "Normal"
$ws.Range("E42").Value2 = "  -1.52%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value2 = "123.47"
$c.Style = "Normal"
$ws.Range("E43").Value2 = "  +10.55%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value2 = "22.40"
$c.Style = "Normal"
$ws.Range("E44").Value2 = "  -3.43%  "
$ws.Range("E45").Value2 = "  -3.89%  "
$ws.Range("D46").Value2 = "2.108.89"
$ws.Range("E46").Value2 = "  -3.03%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value2 = "3.32"
$c.Style = "Normal"
$ws.Range("E47").Value2 = "  -5.20%  "
$ws.Range("E48").Value2 = "  -8.91%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value2 = "0.237"
$c.Style = "Normal"
$ws.Range("E49").Value2 = "  -5.03%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value2 = "0.0336"
$c.Style = "Normal"
$ws.Range("E50").Value2 = "  -4.13%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value2 = "0.914"
$c.Style = "Normal"
$ws.Range("E51").Value2 = "  -3.34%  "
